$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.393.69"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3
$ws.Range("D3").Value = "1.931.35"
$ws.Range("E3").Value = "  -3.97%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.89"
$ws.Range("E5").Value = "  -2.11%  "

# Row 6
$ws.Range("E6").Value = "  -3.87%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.66"
$ws.Range("E8").Value = "  -9.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.365"
$ws.Range("E9").Value = "  -4.97%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.29"
$ws.Range("E10").Value = "  -3.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("E11").Value = "  +4.28%  "

# Row 12
$ws.Range("E12").Value = "  -0.66%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.814"
$ws.Range("E13").Value = "  -7.64%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.30"
$ws.Range("E14").Value = "  -7.16%  "

# Row 15
$ws.Range("D15").Value = "2.211.90"
$ws.Range("E15").Value = "  -4.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.33"
$ws.Range("E16").Value = "  -5.88%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.19"
$ws.Range("E17").Value = "  -6.28%  "

# Row 18
$ws.Range("D18").Value = "1.977.91"
$ws.Range("E18").Value = "  -1.64%  "

# Row 19
$ws.Range("D19").Value = "36.317.28"
$ws.Range("E19").Value = "  -0.59%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.15"
$ws.Range("E20").Value = "  -3.81%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0860"
$ws.Range("E21").Value = "  -1.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.65"
$ws.Range("E22").Value = "  -4.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").Value = "  -7.08%  "

# Row 24
$ws.Range("E24").Value = "  +0.10%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("E25").Value = "  -3.28%  "

# Row 26
$ws.Range("E26").Value = "  -2.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  -7.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.50"
$ws.Range("E28").Value = "  +2.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.22"
$ws.Range("E29").Value = "  -4.61%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.122"
$ws.Range("E30").Value = "  -7.77%  "

# Row 31
$ws.Range("E31").Value = "  -3.48%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  -3.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.63"
$ws.Range("E33").Value = "  -7.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0623"
$ws.Range("E34").Value = "  -1.50%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.27"
$ws.Range("E35").Value = "  -5.38%  "

# Row 36
$ws.Range("E36").Value = "  +0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.98"
$ws.Range("E37").Value = "  -7.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -2.90%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.13"
$ws.Range("E39").Value = "  -9.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("E40").Value = "  -8.61%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0967"
$ws.Range("E41").Value = "  -3.94%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.85"
$ws.Range("E42").Value = "  -1.91%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.17"
$ws.Range("E43").Value = "  -7.51%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0207"
$ws.Range("E44").Value = "  -4.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.66"
$ws.Range("E45").Value = "  -6.13%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.340.54"
$ws.Range("E46").Value = "  -1.46%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.03"
$ws.Range("E47").Value = "  -8.60%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.08"
$ws.Range("E48").Value = "  -9.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.20"
$ws.Range("E49").Value = "  -6.31%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("E50").Value = "  -2.96%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.52"
$ws.Range("E51").Value = "  +1.88%  "
